$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sprint1")

# Move the 4 task rows that belong to "User Story/ Feature #1" (rows 5-8)
# down into the first task slot of each of the following user-story
# sections (rows 10, 14, 19, 24), restoring rows 5,6,7,8 back to the
# blank "Task" placeholder content used elsewhere on the sheet.

# --- capture source values before they get overwritten ---
$a5 = $ws.Range("A5").Value2
$c5 = $ws.Range("C5").Value2
$d5 = $ws.Range("D5").Value2

$a6 = $ws.Range("A6").Value2
$c6 = $ws.Range("C6").Value2
$d6 = $ws.Range("D6").Value2

$a7 = $ws.Range("A7").Value2
$c7 = $ws.Range("C7").Value2
$d7 = $ws.Range("D7").Value2

$a8 = $ws.Range("A8").Value2
$c8 = $ws.Range("C8").Value2
$d8 = $ws.Range("D8").Value2

# --- write them into their new homes ---
$ws.Range("A10").Value = $a5
$ws.Range("C10").Value = $c5
$ws.Range("D10").Value = $d5

$ws.Range("A14").Value = $a6
$ws.Range("C14").Value = $c6
$ws.Range("D14").Value = $d6

$ws.Range("A19").Value = $a7
$ws.Range("C19").Value = $c7
$ws.Range("D19").Value = $d7

$ws.Range("A24").Value = $a8
$ws.Range("C24").Value = $c8
$ws.Range("D24").Value = $d8

# --- reset the source rows back to the blank "Task" template ---
$ws.Range("A5").Value = "Task"
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = ""

$ws.Range("A6").Value = "Task"
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = ""

$ws.Range("A7").Value = "Task"
$ws.Range("B7").Value = ""
$ws.Range("C7").Value = ""
$ws.Range("D7").Value = ""

$ws.Range("A8").Value = "Task"
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = ""

# B19 no longer exists on the destination row once the data moved in
$ws.Range("B19").Value = ""

# Update the remembered selection/active cell for the sheet
$ws.Range("G25").Select()
